# [2 18] update ex12a.py
# Rename "Sheet1" -> "first_sheet" and relocate the data table from A1:C5
# (with its A-column index dropped) to D4:E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "first_sheet"

# Move the Name/Height table (columns B:C, rows 1:5) down-and-right to D4:E8
$ws.Range("B1:C5").Cut($ws.Range("D4"))

# Drop the now-orphaned index column (A1:A5)
$ws.Range("A1:A5").Clear()
